# Generate Report for Handback
# Update the "last generated" timestamps that get stamped during report
# generation for the handback status workbook.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the first data row.
$wsOverview.Range("G2").Value = "2016-09-02 15:20:49"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
# for the first data row.
$wsZhCn.Range("H2").Value = "2016-09-02 15:20:45"
$wsZhCn.Range("K2").Value = "2016-09-02 15:21:09"

# de-de sheet: "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
# for the first data row.
$wsDeDe.Range("H2").Value = "2016-09-02 15:20:49"
$wsDeDe.Range("K2").Value = "2016-09-02 15:21:20"
